# Update automatico via Actualizar 02-07-2021 12-16-35
# Rolls the "Ultimo" timestamp column (D) down: each block of 14 rows
# takes on the timestamp value that used to belong to the block above it,
# and a fresh timestamp is written for the newest block (rows 2-15).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2:D15").Value = 44234.51146017924
$ws.Range("D16:D29").Value = 44234.49027027778
$ws.Range("D30:D43").Value = 44234.46906927083
